$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.281.76'
$ws.Range("E2").Value = '  +0.14%  '
$ws.Range("D3").Value = '1.857.93'
$ws.Range("E3").Value = '  -0.23%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  +0.17%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7037'
$ws.Range("E5").Value = '  +0.39%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '238.15'
$ws.Range("E6").Value = '  +0.21%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.001'
$ws.Range("E7").Value = '  +0.17%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07858'
$ws.Range("E8").Value = '  +2.21%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3034'
$ws.Range("E9").Value = '  -0.86%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '24.51'
$ws.Range("E10").Value = '  +4.93%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08175'
$ws.Range("E11").Value = '  -0.43%  '
$ws.Range("D12").Value = '1.864.54'
$ws.Range("E12").Value = '  +2.20%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.210'
$ws.Range("E13").Value = '  +1.15%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.7151'
$ws.Range("E14").Value = '  -0.39%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '89.57'
$ws.Range("E15").Value = '  -0.03%  '
$ws.Range("D16").Value = '29.340.80'
$ws.Range("E16").Value = '  +0.33%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '5.798'
$ws.Range("E17").Value = '  +0.65%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000007821'
$ws.Range("E18").Value = '  +1.49%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.21'
$ws.Range("E19").Value = '  -0.99%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '237.56'
$ws.Range("E20").Value = '  +0.13%  '
$ws.Range("B21").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C21").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D21").Value = '2.124.15'
$ws.Range("E21").Value = '  +0.71%  '
$ws.Range("B22").Value = 'Dai'
$ws.Range("C22").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.000'
$ws.Range("E22").Value = '  +0.07%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.001'
$ws.Range("E23").Value = '  +0.18%  '
$ws.Range("E24").Value = '  +1.15%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '162.42'
$ws.Range("E25").Value = '  -0.10%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.900'
$ws.Range("E26").Value = '  -1.48%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1419'
$ws.Range("E27").Value = '  -2.56%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.09'
$ws.Range("E28").Value = '  +0.14%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.912'
$ws.Range("E29").Value = '  -2.86%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.396'
$ws.Range("E30").Value = '  -1.26%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.477'
$ws.Range("E31").Value = '  -0.50%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.307'
$ws.Range("E32").Value = '  -3.32%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.048'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05193'
$ws.Range("E34").Value = '  -0.29%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.179'
$ws.Range("E35").Value = '  +1.12%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7120'
$ws.Range("E36").Value = '  +0.22%  '
$ws.Range("E37").Value = '  +0.21%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.676'
$ws.Range("E38").Value = '  +0.73%  '
$ws.Range("E39").Value = '  -0.51%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.690'
$ws.Range("E40").Value = '  -1.74%  '
$ws.Range("D41").Value = '1.140.74'
$ws.Range("E41").Value = '  -1.86%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.9217'
$ws.Range("E42").Value = '  -1.70%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.967'
$ws.Range("E43").Value = '  +1.62%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.4257'
$ws.Range("E44").Value = '  -0.84%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '70.51'
$ws.Range("E45").Value = '  -0.24%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.001'
$ws.Range("E46").Value = '  +0.14%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '101.80'
$ws.Range("E47").Value = '  -1.61%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.5332'
$ws.Range("E48").Value = '  -2.18%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.753'
$ws.Range("E49").Value = '  -3.01%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '9.189'
$ws.Range("E50").Value = '  +0.32%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.000'
$ws.Range("E51").Value = '  +0.08%  '
